# Atualização de bases das ligas, do dia: 29-03-2024 às 17:05
# Swap the data (columns B:AC) between row 115 and row 116,
# keeping column A (the "id" values 113/114) unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$range1 = $ws.Range("B115:AC115")
$range2 = $ws.Range("B116:AC116")

$values1 = $range1.Value2
$values2 = $range2.Value2

$range1.Value2 = $values2
$range2.Value2 = $values1
